# ------------------------------------------------------------------
# Adds player_rating / description / comparison_to_real_players columns
# (G, H, I) to the player_info sheet, and nudges the remembered
# selections on player_info / club_info to match where the author
# left the cursor after entering the new data.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("player_info")
$ws.Activate()

# --- Header row (row 1) ---
$ws.Range("G1").Value = "player_rating"
$ws.Range("H1").Value = "description"
$ws.Range("I1").Value = "comparison_to_real_players"

# --- Per-player rating / description / comparison (rows 2-15) ---
$ws.Range("G2").Value = 88
$ws.Range("H2").Value = "A midfielder with a forward's flair, Khalis is known for his rocket shots and feet that move faster than a squirrel on espresso."
$ws.Range("I2").Value = "Kevin De Bruyne, Heung Min Son"

$ws.Range("G3").Value = 90
$ws.Range("H3").Value = "The calm and collected Faris is a technical wizard in defense, with long balls that could reach the moon and back."
$ws.Range("I3").Value = "Virgil van Dijk, Toni Kroos"

$ws.Range("G4").Value = 81
$ws.Range("H4").Value = "Hafiz is the reliable goal machine—he scores goals like it’s his job (because it is), and his defensive work rate is top-notch too!"
$ws.Range("I4").Value = "Roberto Firmino, N'Golo Kanté"

$ws.Range("G5").Value = 83
$ws.Range("H5").Value = "Defensive-minded and as strong as an ox, Danish also doubles as a goalkeeper when needed—talk about versatility!"
$ws.Range("I5").Value = "John Terry, Peter Schmeichel"

$ws.Range("G6").Value = 82
$ws.Range("H6").Value = "With flair, speed, and the ability to win headers like a giant, Imran is the forward that defenders dread facing."
$ws.Range("I6").Value = "Kai Havertz, Cristiano Ronaldo"

$ws.Range("G7").Value = 79
$ws.Range("H7").Value = "Isa's stamina could put an Energizer Bunny to shame, and his strength makes him a brick wall on the field."
$ws.Range("I7").Value = "Kyle Walker, N’Golo Kanté"

$ws.Range("G8").Value = 85
$ws.Range("H8").Value = "Muk is clinical in front of goal and technical with his feet, making him the kind of forward every team dreams of having."
$ws.Range("I8").Value = "Harry Kane, Riyad Mahrez"

$ws.Range("G9").Value = 86
$ws.Range("H9").Value = "Abdullah’s technical skills and calm demeanor in goal make him the goalkeeper who could stay cool even in a snowstorm."
$ws.Range("I9").Value = "Ederson, Xabi Alonso"

$ws.Range("G10").Value = 85
$ws.Range("H10").Value = "Farhan is a shot-stopper extraordinaire—he saves shots like he's got a magnetic glove."
$ws.Range("I10").Value = "Jan Oblak, Keylor Navas"

$ws.Range("G11").Value = 81
$ws.Range("H11").Value = "Ashraf is a defensive powerhouse, strong enough to stop aerial threats and still have time to grab a snack."
$ws.Range("I11").Value = "Sergio Ramos, Giorgio Chiellini"

$ws.Range("G12").Value = 82
$ws.Range("H12").Value = "Hamizan's positioning is so good, he could find the perfect spot in a crowded room, and his passing is like a GPS system for his teammates."
$ws.Range("I12").Value = "Thomas Müller, David Silva"

$ws.Range("G13").Value = 87
$ws.Range("H13").Value = "With his technical skills and speed, Hanif plays with flair that makes defenders look like they're stuck in quicksand."
$ws.Range("I13").Value = "Eden Hazard, Raheem Sterling"

$ws.Range("G14").Value = 81
$ws.Range("H14").Value = "Nabil is a defensive rock, strong and dominant, and he’s the type of player who makes attackers rethink their life choices."
$ws.Range("I14").Value = "Mats Hummels, Virgil van Dijk"

$ws.Range("G15").Value = 79
$ws.Range("H15").Value = "Gan is the workhorse of the team—strong, hardworking, and always ready to put in the extra effort, even if it means running through a wall."
$ws.Range("I15").Value = "James Milner, Claude Makélélé"

# --- Column widths: best-fit the new columns to their content, ---
# --- matching the widths Excel computed when the data was entered. ---
$ws.Columns.Item(4).ColumnWidth = 6.333333333333333
$ws.Columns.Item(7).ColumnWidth = 10.333333333333334
$ws.Columns.Item(8).ColumnWidth = 111.33333333333333
$ws.Columns.Item(9).ColumnWidth = 26

# --- Restore the cursor position left on player_info after data entry ---
$ws.Range("C17").Select()

# --- club_info sheet cursor was also moved while reviewing the update ---
$wsClub = $wb.Worksheets.Item("club_info")
$wsClub.Activate()
$wsClub.Range("H24").Select()

# --- Leave player_info as the active/visible sheet (tab selected) ---
$ws.Activate()

Write-Output "player_info: ratings/descriptions/comparisons added"
